$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values for columns D, L, M, N, O, P, R, S for rows 2-41
$cols = @("D","L","M","N","O","P","R","S")
$orig = @{}
for ($r = 2; $r -le 41; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $orig[$r] = $rowVals
}

# Mapping: destination row -> source row (based on original data)
$mapping = @{
    2 = 38
    3 = 21
    4 = 17
    5 = 15
    6 = 14
    7 = 11
    8 = 4
    9 = 36
    10 = 31
    11 = 34
    12 = 20
    13 = 9
    14 = 41
    15 = 27
    16 = 6
    17 = 40
    18 = 2
    19 = 35
    20 = 5
    21 = 28
    22 = 30
    23 = 19
    24 = 25
    25 = 12
    26 = 32
    27 = 22
    28 = 26
    29 = 39
    30 = 10
    31 = 24
    32 = 7
    33 = 23
    34 = 16
    35 = 18
    36 = 8
    37 = 13
    38 = 29
    39 = 3
    40 = 37
    41 = 33
}

foreach ($dest in $mapping.Keys) {
    $src = $mapping[$dest]
    $srcVals = $orig[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$dest").Value = $srcVals[$c]
    }
}
